$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "<rdkit.Chem.rdchem.Mol object at 0x7f61ac80f4c0>"
$ws.Range("D3").Value = "<rdkit.Chem.rdchem.Mol object at 0x7f61ac80f610>"
$ws.Range("D4").Value = "<rdkit.Chem.rdchem.Mol object at 0x7f61ac80f5a0>"
$ws.Range("D5").Value = "<rdkit.Chem.rdchem.Mol object at 0x7f61ac80f680>"
$ws.Range("D6").Value = "<rdkit.Chem.rdchem.Mol object at 0x7f61ac80f6f0>"
